$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the EMAIL column: remove the header label and the per-row mailto
# hyperlinks (and their display text) that were built from sample data.
$ws.Hyperlinks.Delete()
$ws.Range("M1:M6").ClearContents()

# The EMAIL column is no longer cramped by mailto text, so widen it.
$ws.Columns.Item(13).ColumnWidth = 27.3

# Reposition the saved selection/cursor, as recorded by the author.
$ws.Range("L10").Select() | Out-Null
